# Actualización automática 2025-11-27 15:30:08
# Registers a new sale of 1950.92 for client "HERRERA CAICEDO LUIS FRANKLIN"
# (advisor ALMEIDA CUATIN JHONATHANN CARLOS) under the "240X80 PORCELANATO"
# group, and propagates the resulting totals across the three sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" -------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

# New sale amount for this client in the "240X80 PORCELANATO" column (D16)
$wsGrupo.Range("D16").Value = 1950.92

# Update the "X de 36" counter in the totals row for column D (D38)
$wsGrupo.Range("D38").Value = "2 de 36"

# --- Sheet "VENTA MENSUAL" ----------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

# noviembre column (F) for the same client
$wsMensual.Range("F16").Value = -5874.77

# noviembre total row
$wsMensual.Range("F38").Value = 3665.239999999999

# --- Sheet "CUMPLIMIENTO MENSUAL" ---------------------------------------
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 3: "240X80 PORCELANATO"
$wsCumpl.Range("D3").Value = 2317.26
$wsCumpl.Range("E3").Value = 1275.25
$wsCumpl.Range("F3").Value = 0.6450253443970929

# Row 12: "PORCELANATO"
$wsCumpl.Range("D12").Value = 1238.28
$wsCumpl.Range("E12").Value = 29975.72
$wsCumpl.Range("F12").Value = 0.03967066060101237

# Row 14: "TOTAL"
$wsCumpl.Range("D14").Value = 4635.110000000001
$wsCumpl.Range("E14").Value = 35644.45164865474
$wsCumpl.Range("F14").Value = 0.1150734965894249
